$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 14, pushing existing rows 14-37 down to 16-39
$insertRange = $ws.Range("A14:T15")
$insertRange.EntireRow.Insert()

# Row 14: Carson
$ws.Cells.Item(14, 1).Value = 1
$ws.Cells.Item(14, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(14, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(14, 4).Value = 44579
$ws.Cells.Item(14, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(14, 5).Value = 15
$ws.Cells.Item(14, 6).Value = "Fruta"
$ws.Cells.Item(14, 7).Value = 100103
$ws.Cells.Item(14, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(14, 9).Value = 100103004
$ws.Cells.Item(14, 10).Value = "Durazno"
$ws.Cells.Item(14, 11).Value = "Carson"
$ws.Cells.Item(14, 12).Value = "Segunda"
$ws.Cells.Item(14, 13).Value = 300
$ws.Cells.Item(14, 14).Value = 22000
$ws.Cells.Item(14, 15).Value = 23000
$ws.Cells.Item(14, 16).Value = 22500
$ws.Cells.Item(14, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(14, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(14, 19).Value = 1250
$ws.Cells.Item(14, 20).Value = 18

# Row 15: Rich Lady
$ws.Cells.Item(15, 1).Value = 1
$ws.Cells.Item(15, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(15, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(15, 4).Value = 44579
$ws.Cells.Item(15, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(15, 5).Value = 15
$ws.Cells.Item(15, 6).Value = "Fruta"
$ws.Cells.Item(15, 7).Value = 100103
$ws.Cells.Item(15, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(15, 9).Value = 100103004
$ws.Cells.Item(15, 10).Value = "Durazno"
$ws.Cells.Item(15, 11).Value = "Rich Lady"
$ws.Cells.Item(15, 12).Value = "Segunda"
$ws.Cells.Item(15, 13).Value = 250
$ws.Cells.Item(15, 14).Value = 22000
$ws.Cells.Item(15, 15).Value = 23000
$ws.Cells.Item(15, 16).Value = 22500
$ws.Cells.Item(15, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(15, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(15, 19).Value = 1250
$ws.Cells.Item(15, 20).Value = 18
